$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue "D2" "63.062.77"
Set-TextValue "E2" "  -2.04%  "
Set-TextValue "D3" "3.125.40"
Set-TextValue "E3" "  -0.56%  "
Set-TextValue "E4" "  +0.10%  "
Set-TextValue "D5" "594.78"
Set-TextValue "E5" "  -2.62%  "
Set-TextValue "D6" "136.57"
Set-TextValue "E6" "  -5.09%  "
Set-TextValue "E7" "  +0.05%  "
Set-TextValue "D8" "3.115.95"
Set-TextValue "E8" "  -0.79%  "
Set-TextValue "D9" "0.518"
Set-TextValue "E9" "  -2.46%  "
Set-TextValue "E10" "  -3.65%  "
Set-TextValue "D11" "5.21"
Set-TextValue "E11" "  -4.34%  "
Set-TextValue "E12" "  -3.64%  "
Set-TextValue "E13" "  -3.18%  "
Set-TextValue "D14" "34.18"
Set-TextValue "E14" "  -3.97%  "
Set-TextValue "D15" "3.636.93"
Set-TextValue "E15" "  -0.47%  "
Set-TextValue "E16" "  +1.51%  "
Set-TextValue "D17" "63.032.25"
Set-TextValue "E17" "  -2.00%  "
Set-TextValue "D18" "3.120.76"
Set-TextValue "E18" "  -0.65%  "
Set-TextValue "E19" "  -2.17%  "
Set-TextValue "D20" "476.62"
Set-TextValue "E20" "  -0.26%  "
Set-TextValue "D21" "14.22"
Set-TextValue "E21" "  -3.60%  "
Set-TextValue "E22" "  -3.77%  "
Set-TextValue "D23" "7.68"
Set-TextValue "E23" "  -2.14%  "
Set-TextValue "D24" "87.26"
Set-TextValue "E24" "  +2.25%  "
Set-TextValue "D25" "13.03"
Set-TextValue "E25" "  -5.08%  "
Set-TextValue "D27" "2.73"
Set-TextValue "E27" "  -2.11%  "
Set-TextValue "D28" "7.21"
Set-TextValue "E28" "  -2.74%  "
Set-TextValue "E29" "  -7.79%  "
Set-TextValue "E30" "  -0.90%  "
Set-TextValue "D31" "27.00"
Set-TextValue "E31" "  +0.90%  "
Set-TextValue "E32" "  +0.00%  "
Set-TextValue "E33" "  -8.43%  "
Set-TextValue "E34" "  -4.36%  "
Set-TextValue "E35" "  -3.06%  "
Set-TextValue "E36" "  -2.33%  "
Set-TextValue "D37" "51.96"
Set-TextValue "E37" "  -1.27%  "
Set-TextValue "D38" "0.0₃0713"
Set-TextValue "E38" "  -4.76%  "
Set-TextValue "E39" "  -2.46%  "
Set-TextValue "D40" "422.65"
Set-TextValue "E40" "  -7.41%  "
Set-TextValue "E41" "  -0.83%  "
Set-TextValue "D42" "8.28"
Set-TextValue "E42" "  -0.89%  "
Set-TextValue "B43" "dogwifhat"
Set-TextValue "C43" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D43" "2.68"
Set-TextValue "E43" "  -11.67%  "
Set-TextValue "B44" "Maker"
Set-TextValue "C44" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D44" "2.882.26"
Set-TextValue "E44" "  +0.20%  "
Set-TextValue "D45" "0.267"
Set-TextValue "E45" "  +0.91%  "
Set-TextValue "D46" "2.14"
Set-TextValue "E46" "  -5.53%  "
Set-TextValue "D48" "25.84"
Set-TextValue "E48" "  -2.93%  "
Set-TextValue "E49" "  -0.74%  "
Set-TextValue "D50" "2.29"
Set-TextValue "E50" "  -7.01%  "
Set-TextValue "D51" "118.49"
Set-TextValue "E51" "  -2.07%  "
